$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet is protected; unprotect before editing, then restore protection after
$ws.Unprotect()

# Update the confidential disclaimer date from 2021-03-24 to 2021-03-25
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-25 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2502841431360442
$ws.Range("E2").Value = 0.01448140900195716

$ws.Range("D3").Value = 0.2477722258533101
$ws.Range("E3").Value = 0.01682692307692313

$ws.Range("D4").Value = 0.2506449823862091
$ws.Range("E4").Value = -0.0006179514908080108

$ws.Range("D5").Value = 0.2512986486244365
$ws.Range("E5").Value = -0.004901293396868511

$ws.Range("D6").Value = 0.9999999999999999
$ws.Range("E6").Value = 0.006407136380814915

# Restore sheet protection (sheet was protected before the edit)
$ws.Protect()
